# Refresh the "cryptos" price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values. A few rows also changed rank/position, so their
# Coin name (B) and Link (C) are rewritten too.
#
# Note: several Price values look like plain decimals (e.g. "7.50", "0.999") which
# Excel would otherwise auto-convert to numbers (and silently drop trailing zeros).
# Those are entered with a leading apostrophe so they stay literal text, matching
# the rest of the column (prices with thousands separators, like "68.936.76",
# can never be parsed as numbers anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.936.76"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.927.25"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'603.82"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'167.84"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "3.926.18"
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "'6.47"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").Value = "'37.66"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "4.581.06"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "3.953.99"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "68.935.52"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'7.50"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'17.48"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "'11.10"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "'493.70"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "'0.731"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'0.0000168"
$ws.Range("E24").Value = "  +4.82%  "
$ws.Range("D25").Value = "'84.86"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'12.12"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D30").Value = "'2.97"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "4.076.89"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").Value = "'2.39"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.77"
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'32.03"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").Value = "3.883.92"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'1.04"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'5.97"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'3.27"
$ws.Range("E40").Value = "  +7.82%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.319"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "'435.16"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "'2.00"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'48.12"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "'8.60"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'143.35"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000270"
$ws.Range("E49").Value = "  +18.78%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.825.06"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'26.04"
$ws.Range("E51").Value = "  +4.28%  "

Write-Output "Updated cryptos table: 104 cell(s) across rows 2-51."
